$wb = $excel.ActiveWorkbook

# "meta" sheet holds key/value settings (tab, type, title, y_title, ...)
$meta = $wb.Worksheets.Item("meta")

# Row 5 was an empty templated row (A5 styled like the other key cells).
# Fill it in with the new "style" / "default" setting, copying the key
# style (bold/orange) from an existing key cell (A4) onto A5.
$meta.Range("A4").Copy()
$meta.Range("A5").PasteSpecial(-4122)  # xlPasteFormats
$meta.Range("A5").Value = "style"
$meta.Range("B5").Value = "default"

# Add a new empty templated row below, matching the previous look of A5
# (same style, no value) so future settings can be appended there.
$meta.Range("A5").Copy()
$meta.Range("A6").PasteSpecial(-4122)  # xlPasteFormats
$meta.Range("A6").Value = ""

$excel.CutCopyMode = 0
